$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''56.702.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").Value = '''2.325.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''521.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.87%  '
$ws.Range("D6").Value = '''134.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.96%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("D9").Value = '''2.351.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("D10").Value = '''0.105'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.74%  '
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").Value = '''5.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.51%  '
$ws.Range("D13").Value = '''0.343'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").Value = '''23.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '''2.744.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '''56.771.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("D18").Value = '''2.355.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("D19").Value = '''10.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '''324.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.59%  '
$ws.Range("D22").Value = '''6.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E25").Value = '  +8.90%  '
$ws.Range("D26").Value = '''0.995'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '''7.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.26%  '
$ws.Range("D28").Value = '''1.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.35%  '
$ws.Range("D29").Value = '''0.0₃0749'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.74%  '
$ws.Range("E30").Value = '  +5.89%  '
$ws.Range("D31").Value = '''169.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  +1.35%  '
$ws.Range("D37").Value = '''0.926'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("E38").Value = '  +3.86%  '
$ws.Range("E39").Value = '  +8.15%  '
$ws.Range("D40").Value = '''37.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.33%  '
$ws.Range("E41").Value = '  +0.51%  '
$ws.Range("E42").Value = '  +4.66%  '
$ws.Range("D43").Value = '''137.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").Value = '''5.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.40%  '
$ws.Range("D45").Value = '''278.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.59%  '
$ws.Range("E46").Value = '  +2.69%  '
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("E48").Value = '  +2.30%  '
$ws.Range("D49").Value = '''0.0217'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.56%  '
$ws.Range("D50").Value = '''17.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.51%  '
$ws.Range("E51").Value = '  +0.54%  '
